# Updated cryptos list: refresh prices/volume deltas scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.958.11'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.818.06'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'311.90"
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').Value = "'1.006"
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = "'0.4301"
$ws.Range('E7').Value = '  +1.32%  '
$ws.Range('D8').Value = "'0.3697"
$ws.Range('E8').Value = '  +2.39%  '
$ws.Range('D9').Value = "'0.07262"
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').Value = "'0.8653"
$ws.Range('E10').Value = '  +2.51%  '
$ws.Range('D11').Value = '2.029.44'
$ws.Range('E11').Value = '  +13.75%  '
$ws.Range('D12').Value = "'21.03"
$ws.Range('E12').Value = '  +3.83%  '
$ws.Range('D13').Value = "'6.647"
$ws.Range('E13').Value = '  +4.06%  '
$ws.Range('D14').Value = "'5.391"
$ws.Range('E14').Value = '  +2.10%  '
$ws.Range('D15').Value = "'0.06928"
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('D16').Value = "'80.76"
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = "'1.006"
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = "'0.000008870"
$ws.Range('E18').Value = '  +1.68%  '
$ws.Range('D19').Value = "'1.006"
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = "'15.21"
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('D21').Value = '27.005.09'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = "'5.193"
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('D23').Value = "'11.11"
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('D24').Value = '2.250.14'
$ws.Range('E24').Value = '  +12.04%  '
$ws.Range('D25').Value = "'154.03"
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').Value = "'1.885"
$ws.Range('E26').Value = '  -3.70%  '
$ws.Range('D27').Value = "'18.28"
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').Value = "'5.221"
$ws.Range('E28').Value = '  +3.96%  '
$ws.Range('D29').Value = "'1.896"
$ws.Range('E29').Value = '  +17.12%  '
$ws.Range('D30').Value = "'115.18"
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').Value = "'0.08972"
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').Value = "'0.7452"
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('D33').Value = "'1.166"
$ws.Range('E33').Value = '  +6.69%  '
$ws.Range('D34').Value = "'4.415"
$ws.Range('E34').Value = '  +1.67%  '
$ws.Range('D35').Value = "'2.803"
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('D36').Value = "'1.010"
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').Value = "'1.126"
$ws.Range('E37').Value = '  +3.68%  '
$ws.Range('D38').Value = "'0.05217"
$ws.Range('E38').Value = '  +2.16%  '
$ws.Range('D39').Value = "'0.01928"
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('D40').Value = "'0.5077"
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('D41').Value = "'2.740"
$ws.Range('E41').Value = '  +6.13%  '
$ws.Range('D42').Value = "'0.1646"
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('D43').Value = "'6.432"
$ws.Range('E43').Value = '  +7.57%  '
$ws.Range('D44').Value = "'8.237"
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = "'107.11"
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'10.38"
$ws.Range('E46').Value = '  +1.93%  '
$ws.Range('D47').Value = "'1.006"
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').Value = "'1.660"
$ws.Range('E48').Value = '  +4.48%  '
$ws.Range('D49').Value = "'0.06304"
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = "'0.4568"
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').Value = "'1.818"
$ws.Range('E51').Value = '  +5.96%  '
